$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1947743467933492
$ws.Range("C2").Value = 0.5676959619952494
$ws.Range("J2").Value = 0.02137767220902613
$ws.Range("P2").Value = 0.1258907363420428
$ws.Range("S2").Value = 0.09026128266033254

# Row 3
$ws.Range("C3").Value = 0.04
$ws.Range("J3").Value = 0.016
$ws.Range("P3").Value = 0.796
$ws.Range("S3").Value = 0.148

# Row 4
$ws.Range("P4").Value = 0.6585365853658537
$ws.Range("S4").Value = 0.3414634146341464

# Row 6
$ws.Range("B6").Value = 0.07722007722007722
$ws.Range("D6").Value = 0.003861003861003861
$ws.Range("F6").Value = 0.07335907335907337
$ws.Range("J6").Value = 0.3011583011583012
$ws.Range("O6").Value = 0.01158301158301158
$ws.Range("Q6").Value = 0.1814671814671815
$ws.Range("R6").Value = 0.0694980694980695
$ws.Range("S6").Value = 0.2818532818532818

# Row 7
$ws.Range("B7").Value = 0.1596244131455399
$ws.Range("D7").Value = 0.009389671361502348
$ws.Range("F7").Value = 0.05164319248826291
$ws.Range("J7").Value = 0.09389671361502347
$ws.Range("O7").Value = 0.03286384976525822
$ws.Range("Q7").Value = 0.1784037558685446
$ws.Range("R7").Value = 0.1220657276995305
$ws.Range("S7").Value = 0.352112676056338

# Row 8
$ws.Range("B8").Value = 0.1204081632653061
$ws.Range("D8").Value = 0.02040816326530612
$ws.Range("E8").Value = 0.004081632653061225
$ws.Range("F8").Value = 0.07551020408163266
$ws.Range("J8").Value = 0.1183673469387755
$ws.Range("O8").Value = 0.02244897959183673
$ws.Range("Q8").Value = 0.1775510204081633
$ws.Range("R8").Value = 0.07959183673469387
$ws.Range("S8").Value = 0.3816326530612245

# Row 9
$ws.Range("B9").Value = 0.1421319796954315
$ws.Range("D9").Value = 0.005076142131979695
$ws.Range("F9").Value = 0.07614213197969544
$ws.Range("J9").Value = 0.09137055837563451
$ws.Range("O9").Value = 0.01522842639593909
$ws.Range("Q9").Value = 0.1979695431472081
$ws.Range("R9").Value = 0.1116751269035533
$ws.Range("S9").Value = 0.3604060913705584

# Row 10
$ws.Range("B10").Value = 0.1406480117820324
$ws.Range("D10").Value = 0.02209131075110457
$ws.Range("E10").Value = 0.001472754050073638
$ws.Range("F10").Value = 0.05743740795287187
$ws.Range("J10").Value = 0.1053019145802651
$ws.Range("O10").Value = 0.0235640648011782
$ws.Range("Q10").Value = 0.203240058910162
$ws.Range("R10").Value = 0.09131075110456553
$ws.Range("S10").Value = 0.3549337260677467

# Row 11
$ws.Range("G11").Value = 0.1467065868263473
$ws.Range("J11").Value = 0.09580838323353294
$ws.Range("K11").Value = 0.2005988023952096
$ws.Range("L11").Value = 0.5449101796407185
$ws.Range("S11").Value = 0.01197604790419162

# Row 12
$ws.Range("G12").Value = 0.7553191489361702
$ws.Range("J12").Value = 0.1702127659574468
$ws.Range("K12").Value = 0.02127659574468085
$ws.Range("L12").Value = 0.02659574468085106
$ws.Range("S12").Value = 0.02659574468085106

# Row 13
$ws.Range("G13").Value = 0.6458333333333334
$ws.Range("S13").Value = 0.02083333333333333

# Row 15
$ws.Range("F15").Value = 0.02489626556016597
$ws.Range("H15").Value = 0.1535269709543569
$ws.Range("I15").Value = 0.05809128630705394
$ws.Range("J15").Value = 0.3319502074688797
$ws.Range("K15").Value = 0.06639004149377593
$ws.Range("M15").Value = 0.008298755186721992
$ws.Range("O15").Value = 0.04564315352697095
$ws.Range("S15").Value = 0.3112033195020747

# Row 16
$ws.Range("F16").Value = 0.02641509433962264
$ws.Range("H16").Value = 0.1283018867924528
$ws.Range("I16").Value = 0.1094339622641509
$ws.Range("J16").Value = 0.4679245283018868
$ws.Range("K16").Value = 0.09433962264150944
$ws.Range("M16").Value = 0.003773584905660377
$ws.Range("N16").Value = 0.003773584905660377
$ws.Range("O16").Value = 0.04905660377358491
$ws.Range("S16").Value = 0.1169811320754717

# Row 17
$ws.Range("F17").Value = 0.03105590062111801
$ws.Range("H17").Value = 0.1780538302277433
$ws.Range("I17").Value = 0.07660455486542443
$ws.Range("J17").Value = 0.3768115942028986
$ws.Range("K17").Value = 0.1325051759834369
$ws.Range("M17").Value = 0.03312629399585922
$ws.Range("O17").Value = 0.04968944099378882
$ws.Range("S17").Value = 0.1221532091097308

# Row 18
$ws.Range("F18").Value = 0.02183406113537118
$ws.Range("H18").Value = 0.1572052401746725
$ws.Range("I18").Value = 0.06550218340611354
$ws.Range("J18").Value = 0.462882096069869
$ws.Range("K18").Value = 0.07860262008733625
$ws.Range("M18").Value = 0.004366812227074236
$ws.Range("O18").Value = 0.06550218340611354
$ws.Range("S18").Value = 0.1441048034934498

# Row 19
$ws.Range("F19").Value = 0.02384500745156483
$ws.Range("H19").Value = 0.2250372578241431
$ws.Range("I19").Value = 0.07749627421758569
$ws.Range("J19").Value = 0.3524590163934426
$ws.Range("K19").Value = 0.09836065573770492
$ws.Range("M19").Value = 0.02235469448584203
$ws.Range("N19").Value = 0.0007451564828614009
$ws.Range("O19").Value = 0.07451564828614009
$ws.Range("S19").Value = 0.1251862891207154

Write-Host "Updated team-specific transition matrix values"